$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells before writing literal
# numeric-looking strings, so Excel keeps them as text exactly
# as authored (no auto-number conversion / trailing-zero loss).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "45.393.33"
$ws.Range("E2").Value = "  +7.21%  "
$ws.Range("D3").Value = "2.385.72"
$ws.Range("E3").Value = "  +4.81%  "
$ws.Range("E4").Value = "  +0.52%  "
$ws.Range("D5").Value = "112.02"
$ws.Range("E5").Value = "  +9.01%  "
$ws.Range("D6").Value = "318.00"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("D7").Value = "0.637"
$ws.Range("E7").Value = "  +4.25%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.630"
$ws.Range("E9").Value = "  +5.66%  "
$ws.Range("D10").Value = "42.26"
$ws.Range("E10").Value = "  +10.24%  "
$ws.Range("E11").Value = "  +4.14%  "
$ws.Range("E12").Value = "  +6.44%  "
$ws.Range("E13").Value = "  +5.80%  "
$ws.Range("E14").Value = "  +1.96%  "
$ws.Range("D15").Value = "15.84"
$ws.Range("E15").Value = "  +5.95%  "
$ws.Range("D16").Value = "2.745.85"
$ws.Range("E16").Value = "  +4.84%  "
$ws.Range("D17").Value = "2.381.77"
$ws.Range("E17").Value = "  +4.52%  "
$ws.Range("D18").Value = "45.398.26"
$ws.Range("E18").Value = "  +7.37%  "
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +7.02%  "
$ws.Range("D20").Value = "0.0000108"
$ws.Range("E20").Value = "  +4.89%  "
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "75.15"
$ws.Range("E22").Value = "  +3.68%  "
$ws.Range("D23").Value = "3.54"
$ws.Range("E23").Value = "  +5.31%  "
$ws.Range("D24").Value = "269.99"
$ws.Range("E24").Value = "  +3.33%  "
$ws.Range("E25").Value = "  +8.44%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").Value = "11.28"
$ws.Range("E27").Value = "  +6.77%  "
$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  +10.96%  "
$ws.Range("E29").Value = "  +0.52%  "
$ws.Range("D30").Value = "22.96"
$ws.Range("E30").Value = "  +4.20%  "
$ws.Range("D31").Value = "38.60"
$ws.Range("E31").Value = "  +9.00%  "
$ws.Range("D32").Value = "0.0943"
$ws.Range("E32").Value = "  +11.50%  "
$ws.Range("D33").Value = "169.89"
$ws.Range("E34").Value = "  +17.89%  "
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  +3.74%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "0.118"
$ws.Range("E36").Value = "  +6.97%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  +9.66%  "
$ws.Range("D38").Value = "3.10"
$ws.Range("E38").Value = "  +14.42%  "
$ws.Range("E39").Value = "  +5.97%  "
$ws.Range("D40").Value = "3.96"
$ws.Range("E40").Value = "  +9.03%  "
$ws.Range("D41").Value = "1.75"
$ws.Range("E41").Value = "  +13.54%  "
$ws.Range("D42").Value = "105.19"
$ws.Range("E42").Value = "  +6.69%  "
$ws.Range("E43").Value = "  +7.75%  "
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  +14.47%  "
$ws.Range("D45").Value = "71.21"
$ws.Range("E45").Value = "  +5.04%  "
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").Value = "117.94"
$ws.Range("E47").Value = "  +8.10%  "
$ws.Range("D48").Value = "5.77"
$ws.Range("E48").Value = "  +13.63%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.36"
$ws.Range("E49").Value = "  +9.32%  "
$ws.Range("B50").Value = "MinaProtocolToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/J7st_qGwz+minaprotocoltoken-mina"
$ws.Range("D50").Value = "1.62"
$ws.Range("E50").Value = "  +19.86%  "
$ws.Range("D51").Value = "79.04"
$ws.Range("E51").Value = "  +4.79%  "

# Restore default style so the text-format tweak above leaves
# no visible trace on the saved styles (cells were unstyled).
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
